# Updates cryptos list prices / 1h volume percentages (and a few swapped
# coin rows) in the "cryptos" worksheet to match the refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some "Price" values (column D) are plain decimal-looking strings
# (e.g. "1.00", "204.02") that must stay as literal text, not be coerced
# into numbers by Excel (which would drop formatting / trailing zeros).
# Briefly flip the cell to Text format while assigning, then restore the
# default "Normal" style so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "68.934.87"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "3.638.33"
$ws.Range("E3").Value = "  +2.64%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "204.02"
$ws.Range("E5").Value = "  +9.01%  "
Set-TextValue "D6" "571.28"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "3.632.12"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("E9").Value = "  -0.16%  "
Set-TextValue "D10" "0.687"
$ws.Range("E10").Value = "  +3.56%  "
Set-TextValue "D11" "61.80"
$ws.Range("E11").Value = "  +17.46%  "
$ws.Range("E12").Value = "  +5.89%  "
Set-TextValue "D13" "0.0000291"
$ws.Range("E13").Value = "  +13.67%  "
Set-TextValue "D14" "10.19"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").Value = "4.217.87"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "3.630.73"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("E17").Value = "  +1.13%  "
Set-TextValue "D18" "19.10"
$ws.Range("E18").Value = "  +5.13%  "
$ws.Range("D19").Value = "68.710.65"
$ws.Range("E19").Value = "  +3.90%  "
Set-TextValue "D20" "12.51"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("E21").Value = "  +3.46%  "
Set-TextValue "D22" "408.13"
$ws.Range("E22").Value = "  +4.43%  "
Set-TextValue "D23" "13.15"
$ws.Range("E23").Value = "  +19.75%  "
Set-TextValue "D24" "4.23"
$ws.Range("E24").Value = "  -1.33%  "
Set-TextValue "D25" "86.12"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +3.42%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D27" "4.00"
$ws.Range("E27").Value = "  +14.36%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D28" "12.76"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D30" "9.43"
$ws.Range("E30").Value = "  +6.77%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D31" "8.10"
$ws.Range("E31").Value = "  +14.58%  "
Set-TextValue "D32" "31.89"
$ws.Range("E32").Value = "  +3.59%  "
Set-TextValue "D33" "670.67"
$ws.Range("E33").Value = "  +6.87%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  +1.00%  "
Set-TextValue "D37" "42.55"
$ws.Range("E37").Value = "  +3.65%  "
Set-TextValue "D38" "0.423"
$ws.Range("E38").Value = "  +7.65%  "
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("E40").Value = "  +0.09%  "
Set-TextValue "D41" "3.27"
$ws.Range("E41").Value = "  +17.50%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.137"
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.227.53"
$ws.Range("E43").Value = "  +8.65%  "
Set-TextValue "D44" "2.77"
$ws.Range("E44").Value = "  +11.81%  "
Set-TextValue "D45" "0.998"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  +26.36%  "
Set-TextValue "D47" "2.89"
$ws.Range("E47").Value = "  +16.00%  "
$ws.Range("E48").Value = "  +4.78%  "
$ws.Range("E49").Value = "  +7.04%  "
$ws.Range("E50").Value = "  +2.53%  "
Set-TextValue "D51" "3.09"
